$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.857.18"
$ws.Range("E2").Value = "  -0.91%  "

$ws.Range("D3").Value = "1.891.08"
$ws.Range("E3").Value = "  -1.10%  "

$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "'0.7752"
$ws.Range("E5").Value = "  -3.56%  "

$ws.Range("D6").Value = "'244.67"
$ws.Range("E6").Value = "  +0.40%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").Value = "'0.3142"
$ws.Range("E8").Value = "  -2.86%  "

$ws.Range("D9").Value = "'0.07377"
$ws.Range("E9").Value = "  +2.55%  "

$ws.Range("D10").Value = "'25.30"
$ws.Range("E10").Value = "  -5.19%  "

$ws.Range("D11").Value = "'0.08138"
$ws.Range("E11").Value = "  +0.81%  "

$ws.Range("D12").Value = "'0.7671"
$ws.Range("E12").Value = "  -1.71%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.945.63"
$ws.Range("E13").Value = "  +1.65%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.486"
$ws.Range("E14").Value = "  +2.04%  "

$ws.Range("D15").Value = "'92.26"
$ws.Range("E15").Value = "  -1.25%  "

$ws.Range("D16").Value = "'6.184"
$ws.Range("E16").Value = "  +2.54%  "

$ws.Range("D17").Value = "29.850.29"
$ws.Range("E17").Value = "  -1.04%  "

$ws.Range("D18").Value = "'13.96"
$ws.Range("E18").Value = "  -1.21%  "

$ws.Range("D19").Value = "'245.02"
$ws.Range("E19").Value = "  -1.13%  "

$ws.Range("D20").Value = "'0.000007837"
$ws.Range("E20").Value = "  +0.39%  "

$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("D22").Value = "'8.132"
$ws.Range("E22").Value = "  -1.12%  "

$ws.Range("D23").Value = "2.135.79"
$ws.Range("E23").Value = "  -1.24%  "

$ws.Range("E24").Value = "  -0.16%  "

$ws.Range("D25").Value = "'0.1583"
$ws.Range("E25").Value = "  -3.13%  "

$ws.Range("D26").Value = "'9.421"
$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("D27").Value = "'162.42"
$ws.Range("E27").Value = "  -2.59%  "

$ws.Range("D28").Value = "'18.82"
$ws.Range("E28").Value = "  -0.54%  "

$ws.Range("D29").Value = "'2.044"
$ws.Range("E29").Value = "  -4.19%  "

$ws.Range("D30").Value = "'1.454"
$ws.Range("E30").Value = "  +5.08%  "

$ws.Range("D31").Value = "'1.550"
$ws.Range("E31").Value = "  +0.19%  "

$ws.Range("D32").Value = "'4.489"
$ws.Range("E32").Value = "  -0.68%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.05600"
$ws.Range("E33").Value = "  -2.28%  "

$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'4.091"
$ws.Range("E34").Value = "  -1.16%  "

$ws.Range("E35").Value = "  -2.86%  "

$ws.Range("D36").Value = "'0.7633"
$ws.Range("E36").Value = "  +2.37%  "

$ws.Range("D37").Value = "'1.004"
$ws.Range("E37").Value = "  +0.54%  "

$ws.Range("D38").Value = "'2.645"
$ws.Range("E38").Value = "  -3.12%  "

$ws.Range("D39").Value = "'0.01927"
$ws.Range("E39").Value = "  -1.28%  "

$ws.Range("D40").Value = "'2.786"
$ws.Range("E40").Value = "  -0.85%  "

$ws.Range("D41").Value = "1.164.46"
$ws.Range("E41").Value = "  +12.35%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'74.41"
$ws.Range("E42").Value = "  +1.43%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.4467"
$ws.Range("E43").Value = "  -0.58%  "

$ws.Range("D44").Value = "'5.985"
$ws.Range("E44").Value = "  -0.13%  "

$ws.Range("D45").Value = "'0.8525"
$ws.Range("E45").Value = "  -0.03%  "

$ws.Range("E46").Value = "  -0.04%  "

$ws.Range("D47").Value = "'1.902"
$ws.Range("E47").Value = "  -0.52%  "

$ws.Range("D48").Value = "'102.11"
$ws.Range("E48").Value = "  -0.46%  "

$ws.Range("D49").Value = "'9.924"
$ws.Range("E49").Value = "  -0.46%  "

$ws.Range("D50").Value = "'3.092"
$ws.Range("E50").Value = "  +0.74%  "

$ws.Range("D51").Value = "'7.533"
$ws.Range("E51").Value = "  -0.55%  "
